$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, shifting existing rows (and old row 2..25) down to 3..26
$ws.Rows("2:2").Insert()

# The inserted row inherits formatting from the row above (the bold header row);
# clear that so the new data row looks like the other plain data rows.
$ws.Range("A2:R2").ClearFormats()

# Re-apply the date number-format style used by every other cell in column D
$ws.Range("D3").Copy()
$ws.Range("D2").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the new row 2 data
$ws.Range("A2").Value = 7
$ws.Range("B2").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C2").Value = "Ñuble"
$ws.Range("D2").Value = 44630
$ws.Range("E2").Value = 16
$ws.Range("F2").Value = 100112001
$ws.Range("G2").Value = "Berenjena"
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 60
$ws.Range("K2").Value = 9000
$ws.Range("L2").Value = 9500
$ws.Range("M2").Value = 9250
$ws.Range("N2").Value = "$/caja 60 unidades"
$ws.Range("O2").Value = "Región Metropolitana"
$ws.Range("P2").Value = 154
$ws.Range("Q2").Value = 60
$ws.Range("R2").Value = "Hortaliza"
